# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the d6049727-4991-46e2-b550-c0a797672a61 file (row 7) following a new
# handoff xliff generation pass.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-09-06 20:56:38"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-09-06 20:56:33"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H7").Value = "2016-09-06 20:56:38"
